# Apply updated crypto price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.114.49"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "2.225.65"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "291.79"
$ws.Range("E5").Value = "  -0.53%  "

$ws.Range("D6").Value = "87.64"
$ws.Range("E6").Value = "  +1.32%  "

$ws.Range("E7").Value = "  -0.60%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("D10").Value = "30.41"
$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0780"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.34%  "

$ws.Range("E12").Value = "  +3.29%  "

$ws.Range("E13").Value = "  +0.89%  "

$ws.Range("D14").Value = "2.570.22"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").Value = "13.95"
$ws.Range("E15").Value = "  -1.95%  "

$ws.Range("D16").Value = "2.251.92"
$ws.Range("E16").Value = "  +0.69%  "

$ws.Range("E17").Value = "  -0.45%  "

$ws.Range("D18").Value = "40.070.76"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").Value = "0.0₃0887"
$ws.Range("E19").Value = "  -1.06%  "

$ws.Range("E20").Value = "  +7.20%  "

$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").Value = "65.58"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "236.75"

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("E25").Value = "  +1.30%  "

$ws.Range("D26").Value = "1.82"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("D27").Value = "22.68"
$ws.Range("E27").Value = "  -1.49%  "

$ws.Range("E28").Value = "  -1.38%  "

$ws.Range("E29").Value = "  -0.61%  "

$ws.Range("D30").Value = "156.41"
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").Value = "31.76"
$ws.Range("E31").Value = "  -7.13%  "

$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("E33").Value = "  +1.44%  "

$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("E35").Value = "  -1.36%  "

$ws.Range("D36").Value = "2.89"
$ws.Range("E36").Value = "  +6.55%  "

$ws.Range("D38").Value = "15.71"
$ws.Range("E38").Value = "  -5.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0980"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.71%  "

$ws.Range("E40").Value = "  +1.52%  "

$ws.Range("D41").Value = "2.122.30"
$ws.Range("E41").Value = "  +7.83%  "

$ws.Range("D42").Value = "3.87"
$ws.Range("E42").Value = "  +1.76%  "

$ws.Range("E43").Value = "  -1.67%  "

$ws.Range("E44").Value = "  -1.19%  "

$ws.Range("D45").Value = "17.85"
$ws.Range("E45").Value = "  +9.92%  "

$ws.Range("D46").Value = "9.94"
$ws.Range("E46").Value = "  +3.34%  "

$ws.Range("E47").Value = "  +2.61%  "

$ws.Range("D48").Value = "2.435.67"
$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.45%  "

$ws.Range("D51").Value = "69.44"
$ws.Range("E51").Value = "  -2.36%  "
